# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp
# - Update case counts for Rusia, Ucrania/Bolivia, Armenia, Hungria/Malaui/
#   Republica de Yibuti, Lituania and Estonia
# - Ucrania overtakes Bolivia and Hungria overtakes Malaui & Republica de
#   Yibuti in the ranking (sheet is sorted by "Casos totales" desc), so
#   those rows swap which country they describe.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 28 de Agosto de 2020 a las 09:52"

# --- Row 7: Rusia (rank unchanged) -------------------------------------
$ws.Range("B7").Value = 980405
$ws.Range("C7").Value = 4829
$ws.Range("D7").Value = 798466
$ws.Range("E7").Value = 165025
$ws.Range("G7").Value = 110
$ws.Range("H7").Value = 16914

# --- Rows 29-30: Ucrania overtakes Bolivia ------------------------------
$ws.Range("A29").Value = "Ucrania"
$ws.Range("B29").Value = 114497
$ws.Range("C29").Value = 2438
$ws.Range("D29").Value = 55083
$ws.Range("E29").Value = 56963
$ws.Range("G29").Value = 48
$ws.Range("H29").Value = 2451

$ws.Range("A30").Value = "Bolivia"
$ws.Range("B30").Value = 113129
$ws.Range("C30").Value = 1035
$ws.Range("D30").Value = 52521
$ws.Range("E30").Value = 55817
$ws.Range("G30").Value = 65
$ws.Range("H30").Value = 4791

# --- Row 58: Armenia (rank unchanged) -----------------------------------
$ws.Range("B58").Value = 43451
$ws.Range("C58").Value = 181
$ws.Range("D58").Value = 37264
$ws.Range("E58").Value = 5318
$ws.Range("G58").Value = 5
$ws.Range("H58").Value = 869

# --- Rows 107-109: Hungria overtakes Malaui & Republica de Yibuti -------
$ws.Range("A107").Value = "Hungria"
$ws.Range("B107").Value = 5511
$ws.Range("C107").Value = 132
$ws.Range("D107").Value = 3759
$ws.Range("E107").Value = 1138
$ws.Range("H107").Value = 614

$ws.Range("A108").Value = "Malaui"
$ws.Range("B108").Value = 5496
$ws.Range("D108").Value = 3121
$ws.Range("E108").Value = 2202
$ws.Range("H108").Value = 173

$ws.Range("A109").Value = "Republica de Yibuti"
$ws.Range("B109").Value = 5383
$ws.Range("D109").Value = 5307
$ws.Range("E109").Value = 16
$ws.Range("H109").Value = 60

# --- Row 128: Lituania (rank unchanged) ---------------------------------
$ws.Range("B128").Value = 2810
$ws.Range("C128").Value = 48
$ws.Range("D128").Value = 1816
$ws.Range("E128").Value = 908
$ws.Range("G128").Value = 1
$ws.Range("H128").Value = 86

# --- Row 136: Estonia (rank unchanged) ----------------------------------
$ws.Range("B136").Value = 2343
$ws.Range("C136").Value = 18
$ws.Range("D136").Value = 2076
$ws.Range("E136").Value = 203
